# Recon_Summary_Dec_2025_to_Jan_2026.xlsx
# "DQ Changes" sheet (rows 3-9): loans were re-sorted, shuffling which
# Loan ID / Dec-status / Jan-status / UPB combination sits in each row.
# The Jan-status cell (column C) carries a font color that tracks its own
# text ("Current" => dark green, any DPD bucket => dark red) plus a banded
# fill that simply follows the row (odd data-row => no fill, even data-row
# => light peach), independent of content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DQ Changes")

$COLOR_GREEN = 25600      # RGB(0,100,0)   -> dark green text  (Current)
$COLOR_RED   = 192        # RGB(192,0,0)   -> dark red text    (DPD buckets)
$COLOR_PEACH = 14083324   # RGB(252,228,214) -> light peach fill (banded rows)

# NOTE: this COM-interop interpreter does not bind named (-Param value)
# arguments to function parameters reliably, so the helper below is called
# with positional args.
function Set-DqRow {
    param(
        [int]$Row,
        [string]$LoanId,
        [string]$DecStatus,
        [string]$JanStatus,
        [double]$Upb,
        [bool]$Banded
    )

    $ws.Range("A$Row").Value = $LoanId
    $ws.Range("B$Row").Value = $DecStatus

    $cCell = $ws.Range("C$Row")
    $cCell.Value = $JanStatus

    if ($JanStatus -eq "Current") {
        $cCell.Font.Color = $COLOR_GREEN
    } else {
        $cCell.Font.Color = $COLOR_RED
    }

    if ($Banded) {
        $cCell.Interior.Color = $COLOR_PEACH
    } else {
        $cCell.Interior.Pattern = -4142   # xlNone
    }

    $ws.Range("D$Row").Value = $Upb
}

Set-DqRow 3 "MSR100726" "60 DPD"  "Current"  322836.07 $false
Set-DqRow 4 "MSR100542" "Current" "30 DPD"   251963.15 $true
Set-DqRow 5 "MSR100195" "Current" "30 DPD"   217167.13 $false
Set-DqRow 6 "MSR100869" "30 DPD"  "60 DPD"   308607.16 $true
Set-DqRow 7 "MSR100499" "Current" "30 DPD"   383436.85 $false
Set-DqRow 8 "MSR100289" "30 DPD"  "Current"  382060.03 $true
Set-DqRow 9 "MSR100443" "60 DPD"  "90+ DPD"  326446.17 $false

Write-Host "DQ Changes rows 3-9 updated"
